# Update column F (dSF) values for the rows that changed in the repull/push/mean-calc update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -9
    6  = 0
    7  = 0
    17 = 1
    22 = 19
    23 = -1
    24 = -1
    27 = -2
    28 = -3
    29 = -10
    30 = -5
    31 = -10
    34 = -12
    37 = -2
    38 = 1
    41 = 4
    42 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
